$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'network analysis,classification,data queries,overlay analysis'
$ws.Range("B3").Value = ''
$ws.Range("B4").Value = 'overlay analysis,data queries'
$ws.Range("B5").Value = 'network analysis,data queries'
$ws.Range("B6").Value = 'network analysis,classification,data queries,overlay analysis,data queries'
$ws.Range("B7").Value = 'buffer,overlay analysis,data editing,data queries'
$ws.Range("B8").Value = 'classification'
$ws.Range("B9").Value = 'geometry measurement,data queries'
$ws.Range("B10").Value = 'data queries'
$ws.Range("B11").Value = 'buffer,overlay analysis'
$ws.Range("B12").Value = 'data model conversion,overlay analysis,data queries'
$ws.Range("B13").Value = 'buffer,buffer,overlay analysis,overlay analysis'
$ws.Range("B14").Value = 'overlay analysis'
$ws.Range("B15").Value = 'data queries,network analysis,classification,data queries,overlay analysis,data queries,overlay analysis,data queries,overlay analysis'
$ws.Range("B16").Value = 'network analysis,data queries,network analysis,classification,data queries,overlay analysis'
$ws.Range("B17").Value = 'network analysis,classification,data queries,overlay analysis,data queries,overlay analysis,data queries,overlay analysis'
$ws.Range("B18").Value = 'data queries,buffer,overlay analysis'
$ws.Range("B19").Value = 'classification,data queries,data model conversion,overlay analysis'
$ws.Range("B20").Value = 'geometry measurement,data editing,data queries'
$ws.Range("B21").Value = 'network analysis,classification,data queries,overlay analysis,data queries,overlay analysis'
$ws.Range("B22").Value = 'buffer,overlay analysis,data queries'
$ws.Range("B23").Value = 'data queries,buffer,overlay analysis,data queries'
$ws.Range("B24").Value = 'topography'
$ws.Range("B25").Value = 'overlay analysis,data editing,data queries'
$ws.Range("B26").Value = 'data queries,overlay analysis,data editing,data queries'
$ws.Range("B27").Value = 'data queries,data editing'
$ws.Range("B28").Value = 'generalization,geostatistics  '
$ws.Range("B29").Value = 'data queries,generalization,geostatistics  '
$ws.Range("B30").Value = 'data queries,geostatistics  '
$ws.Range("B31").Value = 'geostatistics  '
$ws.Range("B32").Value = 'data queries,overlay analysis,data editing'
$ws.Range("B33").Value = 'data queries,network analysis,data queries'
$ws.Range("B34").Value = 'overlay analysis,geostatistics  '
$ws.Range("B35").Value = 'network analysis'
$ws.Range("B36").Value = 'buffer,overlay analysis,data queries,geometry measurement,data queries'
$ws.Range("B37").Value = 'buffer,buffer,buffer,overlay analysis,overlay analysis,overlay analysis,geometry measurement,data queries,geometry measurement,data queries,geometry measurement,data queries'
$ws.Range("B38").Value = 'geometry measurement,data queries,buffer,overlay analysis,data queries'

$ws.Rows("39:56").Delete()
